$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.449.33'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.923.57'
$ws.Range('E3').Value = '  +3.54%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.94'
$ws.Range('E5').Value = '  +2.56%  '
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4687'
$ws.Range('E7').Value = '  -0.87%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.62'
$ws.Range('E8').Value = '  +3.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2865'
$ws.Range('E9').Value = '  +3.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06924'
$ws.Range('E10').Value = '  +7.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '105.98'
$ws.Range('E11').Value = '  +19.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.26'
$ws.Range('E12').Value = '  -0.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.896.50'
$ws.Range('E13').Value = '  +2.01%  '
$ws.Range('E14').Value = '  +2.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.180'
$ws.Range('E15').Value = '  +3.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.6565'
$ws.Range('E16').Value = '  +3.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '293.89'
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '30.467.86'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.07'
$ws.Range('E19').Value = '  +1.71%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007647'
$ws.Range('E20').Value = '  +3.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.0000'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.151.04'
$ws.Range('E22').Value = '  +3.45%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.260'
$ws.Range('E24').Value = '  +4.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.220'
$ws.Range('E25').Value = '  +3.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.60'
$ws.Range('E26').Value = '  +14.15%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.299'
$ws.Range('E27').Value = '  +1.94%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.59'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.040'
$ws.Range('E29').Value = '  +5.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.1103'
$ws.Range('E30').Value = '  +6.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.364'
$ws.Range('E31').Value = '  +1.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.147'
$ws.Range('E32').Value = '  +1.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.972'
$ws.Range('E33').Value = '  +1.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05089'
$ws.Range('E34').Value = '  +3.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7397'
$ws.Range('E35').Value = '  +3.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.148'
$ws.Range('E36').Value = '  -0.87%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.738'
$ws.Range('E37').Value = '  +1.91%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02025'
$ws.Range('E38').Value = '  +7.19%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.687'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.058'
$ws.Range('E40').Value = '  +2.21%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8752'
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '108.14'
$ws.Range('E42').Value = '  +1.51%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.815'
$ws.Range('E43').Value = '  +5.70%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4233'
$ws.Range('E45').Value = '  +2.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '52.86'
$ws.Range('E46').Value = '  +25.56%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '68.00'
$ws.Range('E47').Value = '  +7.58%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.207'
$ws.Range('E48').Value = '  -1.03%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.202'
$ws.Range('E49').Value = '  +6.30%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1213'
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.72'
$ws.Range('E51').Value = '  +0.85%  '
